# knapsackProblemOnePage.pptx - "Add files via upload" edit
#
# Moves three shapes in the "Resultados/Conclusão" column, grows & rewrites
# the GA-conclusion text box with a second paragraph, and (best-effort)
# enables slide guides at the presentation level.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Conclusão" rectangle (Retângulo 15) moves up -------------------------
$shConclusao = $s.Shapes.Item(9)
$shConclusao.Top = 196.07347106933594

# --- "Trabalhos Citados" rectangle (Retângulo 16) moves down ---------------
$shTrabalhos = $s.Shapes.Item(10)
$shTrabalhos.Top = 446.3885192871094

# --- GA conclusion text box (CaixaDeTexto 26) -------------------------------
# Gets new wording (adds "sempre") plus an additional paragraph, and grows
# downward/taller to make room for it.
$shGA = $s.Shapes.Item(17)
$shGA.TextFrame.TextRange.Text = "Foi possível resolver o problema da mochila e também do cargo de forma eficiente utilizado o algoritmo GA. O valor do fitness sempre se estabiliza num patamar ótimo dependendo dos dados de entrada. Se há uma população inicial baixa ou extremamente alta, há uma grande tendência de uma parada prematura, isso ocorre também com cromossomos com baixo números de alelos.`rO resultado mais significativo deste projeto foi  o entendimento que o algoritmo GA pode ser estendido para  solucionar problemas onde hajam situações onde envolvam otimização de espaços que envolvam múltiplas variáveis, não apenas peso e valor."
$shGA.Top = 223.48764038085938
$shGA.Height = 225.37969970703125

# --- "Referências" text box (CaixaDeTexto 27) moves down --------------------
$shRef = $s.Shapes.Item(18)
$shRef.Top = 478.5649108886719

# --- Presentation-level slide guide list (cosmetic, empty p15:sldGuideLst) --
# Best-effort: some hosts only persist this once a guide collection is
# touched. Harmless if unsupported.
try {
    $ppt.DisplayGuides = $true
    $null = $p.Guides.Add(1, 3)
} catch {
}
